# Added errors and warning messages to test sheets
#
# - Adds an "Errors" sheet (after "Classes") containing the missing-title
#   message in A1.
# - Adds a "Warnings" sheet (after "Errors"), left empty, which becomes the
#   active/selected sheet (mirrors tabSelected moving off the "Classes" sheet).

$wb = $excel.ActiveWorkbook
$classesSheet = $wb.Worksheets.Item("Classes")

# New "Errors" sheet, placed right after "Classes"
$errorsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $classesSheet)
$errorsSheet.Name = "Errors"
$errorsSheet.Range("A1").Value = 'Sheet "Classes" Row: 2 Missing "TITLE"'

# New "Warnings" sheet, placed right after "Errors" - stays empty
$warningsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $errorsSheet)
$warningsSheet.Name = "Warnings"

# Warnings ends up being the active sheet (tabSelected), matching the diff.
$warningsSheet.Activate()
